$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(1)
$sh.Name = $sh.Name
